# Apply updated dataframe values to the active worksheet.
# These correspond to recalculated precision/recall/fmeasure and
# excel/excel_selected/duplicated columns used for the bar chart in
# the accompanying notebook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (state = s0)
$ws.Range("E2").Value = 14.29
$ws.Range("P2").Value = 1
$ws.Range("R2").Value = 7
$ws.Range("U2").Value = 14.29

# Row 3 (state = -)
$ws.Range("P3").Value = 3

# Row 4 (state = union)
$ws.Range("P4").Value = 294
$ws.Range("R4").Value = 1

# Row 5 (state = elcompendex)
$ws.Range("C5").Value = 2.41
$ws.Range("E5").Value = 4.46
$ws.Range("P5").Value = 245
$ws.Range("Q5").Value = 5
$ws.Range("R5").Value = 4
$ws.Range("S5").Value = 2.41
$ws.Range("U5").Value = 4.46

# Row 6 (state = webofscience)
$ws.Range("C6").Value = 2.38
$ws.Range("E6").Value = 4.35
$ws.Range("P6").Value = 210
$ws.Range("Q6").Value = 5
$ws.Range("S6").Value = 2.38
$ws.Range("U6").Value = 4.35

# Row 7 (state = wiley)
$ws.Range("C7").Value = 0.8099999999999999
$ws.Range("E7").Value = 1.39
$ws.Range("P7").Value = 115
$ws.Range("R7").Value = 9
$ws.Range("S7").Value = 0.8099999999999999
$ws.Range("U7").Value = 1.39

# Row 8 (state = sciencedirect)
$ws.Range("C8").Value = 9.09
$ws.Range("E8").Value = 14.43
$ws.Range("P8").Value = 68
$ws.Range("Q8").Value = 7
$ws.Range("R8").Value = 9
$ws.Range("S8").Value = 9.09
$ws.Range("U8").Value = 14.43

# Row 9 (state = acm)
$ws.Range("S9").Value = 1.82
$ws.Range("U9").Value = 3.56
